$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.424225807189941
$ws.Range("B1").Value = 1.605999708175659
$ws.Range("C1").Value = 1.959909677505493
$ws.Range("D1").Value = 2.661046743392944
$ws.Range("E1").Value = 6.697809219360352
